# "vizualizare repartizare din administrator"
# Updates the student/room assignment ("Repartizare") shown on the
# "Informatica" sheet, refreshing the uid -> room/seat mapping for the
# "C1" room block (rows 3-9) and the "C12" room block (rows 17-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Informatica")

# --- "C1" room block (rows 3-9) ---
$ws.Range("D4").Value = "uid023"
$ws.Range("D5").Value = "uid026"
$ws.Range("D6").Value = "uid025"
$ws.Range("D7").Value = "uid022"

$ws.Range("C8").Value = "uid021"
$ws.Range("D8").Value = "uid027"
$ws.Range("E8").Value = "uid029"

$ws.Range("C9").Value = "uid024"
$ws.Range("D9").Value = "uid028"

# --- "C12" room block (rows 17-19) ---
$ws.Range("E17").Value = "uid003"
$ws.Range("F17").Value = "uid006"
$ws.Range("G17").Value = ""

$ws.Range("C18").Value = "uid004"
$ws.Range("D18").Value = "uid005"
$ws.Range("E18").Value = "uid008"
$ws.Range("F18").Value = "uid011"

$ws.Range("C19").Value = "uid009"
$ws.Range("G19").Value = "uid010"
